# tickets_pideu.xlsx update (v.27)
# - insert a blank row 2 (between the header row and the first data row)
# - clear the stray empty H159/I159 cells
# - append new incident rows 160-164 (2024-05-22 paletizer / screw / elevator faults)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force plain text so date/time-looking strings ("2024-05-22", "11:22:43", ...)
    # are not reinterpreted as real Excel date/time serials, matching the rest
    # of the sheet (every existing cell is stored as literal text).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Insert an empty row 2, keeping row 3 (and everything after it) in place ---
$ws.Rows.Item(2).OutlineLevel = 1
$ws.Rows.Item(2).OutlineLevel = 0

# --- Row 159 previously ended with two stray empty cells in H/I; remove them ---
$ws.Range("H159:I159").ClearContents()

# --- New row 160 ---
Set-TextValue "A160" "2024-05-22"
Set-TextValue "B160" "11:22:43"
Set-TextValue "C160" "Fallo en paletizador"
Set-TextValue "D160" "-"
Set-TextValue "E160" "-"
Set-TextValue "F160" "-"
Set-TextValue "G160" "-"
Set-TextValue "H160" "11:22:51"
Set-TextValue "I160" "0:00:08"

# --- New row 161 ---
Set-TextValue "A161" "2024-05-22"
Set-TextValue "B161" "11:22:53"
Set-TextValue "C161" "No atornilla tapa"
Set-TextValue "D161" "-"
Set-TextValue "E161" "-"
Set-TextValue "F161" "-"
Set-TextValue "G161" "-"
Set-TextValue "H161" "11:22:56"
Set-TextValue "I161" "0:00:03"

# --- New row 162 ---
Set-TextValue "A162" "2024-05-22"
Set-TextValue "B162" "11:23:21"
Set-TextValue "C162" "Fallo en elevador"
Set-TextValue "D162" "-"
Set-TextValue "E162" "-"
Set-TextValue "F162" "-"
Set-TextValue "G162" "-"
Set-TextValue "H162" "11:23:35"
Set-TextValue "I162" "0:00:14"

# --- New row 163 ---
Set-TextValue "A163" "2024-05-22"
Set-TextValue "B163" "11:25:01"
Set-TextValue "C163" "Fallo en elevador_3"
Set-TextValue "D163" "-"
Set-TextValue "E163" "-"
Set-TextValue "F163" "-"
Set-TextValue "G163" "-"
Set-TextValue "H163" "11:25:03"
Set-TextValue "I163" "0:00:02"

# --- New row 164 ---
Set-TextValue "A164" "2024-05-22"
Set-TextValue "B164" "11:25:42"
Set-TextValue "C164" "Fallo en elevador_3"
Set-TextValue "D164" "-"
Set-TextValue "E164" "-"
Set-TextValue "F164" "-"
Set-TextValue "G164" "-"
Set-TextValue "H164" "11:25:43"
Set-TextValue "I164" "0:00:01"
